$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OpticalTransmitterBoard")

# --- Row 2: SN75453BD -> SN75454BD OR-logic Driver; its retailer link / RS part number are removed ---
$ws.Range("A2").Value2 = "SN75454BD OR-logic Driver"
$ws.Range("D2").Value2 = ""
$ws.Range("F2").Value2 = ""

# --- Row 3: current limit resistor value 1206 82 Ohm -> 1206 60 Ohm ---
$ws.Range("B3").Value2 = "1206 60 Ohm"

# --- Row 6: SMD Test point hook quantity 36 -> 42 (grounded test points added) ---
$ws.Range("C6").Value2 = 42

# --- Row 10: 0603 100pF formfactor label 100pF 0603 -> 100p ---
$ws.Range("B10").Value2 = "100p"

# --- Row 11: replace "5V, 3A Power supply" with the new decoupling capacitor part ---
$ws.Range("A11").Value2 = "0603 4.7uF"
$ws.Range("B11").Value2 = "capacitor 0603"
$ws.Range("C11").Value2 = 36
$ws.Range("E11").Value2 = "301-98-110"

# --- Hyperlinks: this engine's Hyperlinks.Delete() on a single range clears the whole
# sheet's collection, so drop them all and re-add the ones that should survive,
# dropping the SN75453/4 datasheet link (D2) and repointing D11 at the new capacitor. ---
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D8"), "https://www.elfa.se/sv/kabel-till-kretskort-kopplingsplint-13-31mm-26-16awg-08mm-stiftavstand-poler-rnd-connect-rnd-205-00232/p/30043058?q=&pos=2&origPos=2&origPageSize=50&track=true") | Out-Null
$ws.Range("D8").Style = "Hyperlänk"

$ws.Hyperlinks.Add($ws.Range("D9"), "https://se.rs-online.com/web/p/pcb-sockets/7679647/") | Out-Null
$ws.Range("D9").Style = "Hyperlänk"

$ws.Hyperlinks.Add($ws.Range("D7"), "https://www.elfa.se/sv/keramisk-kondensator-100nf-50vdc-0603-10-rnd-components-rnd-1500603b104k500nt/p/30086429?q=&pos=2&origPos=2&origPageSize=50&track=true") | Out-Null
$ws.Range("D7").Style = "Hyperlänk"

$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.elfa.se/sv/tjockfilmsmotstand-motstand-ytmontering-0603-51ohm-100mw-rnd-components-rnd-1550603saf510jt5e/p/30056626?q=*&pos=9&origPos=2940&origPageSize=50&track=true") | Out-Null
$ws.Range("D5").Style = "Hyperlänk"

$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.elfa.se/sv/ytmonterade-tjockfilmsmotstand-foer-hoegeffekt-1206-82ohm-500mw-rnd-components-rnd-155hp06w2f820jt5e/p/30111989?q=&pos=9&origPos=1641&origPageSize=50&track=true") | Out-Null
$ws.Range("D3").Style = "Hyperlänk"

$ws.Hyperlinks.Add($ws.Range("D10"), "https://www.elfa.se/sv/keramisk-kondensator-100pf-100v-0603-avx-06031a101jat2a/p/30203076?q=kondensator&pos=1&origPos=366&origPageSize=50&track=true") | Out-Null
$ws.Range("D10").Style = "Hyperlänk"

$ws.Range("D11").Value2 = ""
$ws.Hyperlinks.Add($ws.Range("D11"), "https://www.elfa.se/sv/keramisk-kondensator-7uf-10v-0603-10-epcos-c1608x7s1a475k080ac/p/30198110?q=4.7u&pos=20&origPos=20&origPageSize=50&track=true") | Out-Null
$ws.Range("D11").Style = "Hyperlänk"

# --- Selection moved to E15 on the active sheet ---
$ws.Range("E15").Select()
